$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "296.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.74%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.24%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.044"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.67%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07565"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.66%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.398"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.71%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.592"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.73%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.26%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1214"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.83%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1841"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.17%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08996"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.69%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04022"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.74%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1054"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.11%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001292"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.22%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005804"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.64%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.362"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.26%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.927"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.89%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1421"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.92%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.14%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04060"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.30%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.44%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003962"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "4.44%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.87%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.02%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02411"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.60%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05215"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.96%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006258"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-4.29%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007790"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.53%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.66%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007545"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.36%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007848"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.90%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.2974"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.38%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006786"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.66%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04565"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "203.48%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004205"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.07%"
